$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DOI cell (B6) used to store the full "doi.org/..." string both as the
# cell text and as the hyperlink target. After the PubMed/PMC crawl update,
# the cell text becomes the bare DOI while the hyperlink keeps the old
# "doi.org/..." text as its display label.
$ws.Range("B6").Value = "10.1093/cercor/bhh186"

# Rebuild the hyperlink on B6 so it carries an explicit display string
# (or the raw DOI text would otherwise also show up as the link label).
$ws.Hyperlinks.Delete()
$link = $ws.Range("B6").Hyperlinks.Item(1)
$link.Address = "https://doi.org/10.1093/cercor/bhh186"
$link.TextToDisplay = "doi.org/10.1093/cercor/bhh186"

# Leave the selection on the DOI cell, matching where the curator was
# last working.
$ws.Range("B6").Select()
